$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Add the new "BIG O" worksheet, positioned after the current last
# sheet (QUEUES), which matches the target sheet order / activeTab.
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "BIG O"

# Colour used throughout: rgb FF232629 -> OLE BGR integer
$darkColor = 2696739

# Rows 1-6: "O(...)" in the base cell font, followed by
# ' == "word"' rendered in a lighter/larger "Inherit" font run.
$rows = @(
    @("O(1)", " == ""constant"""),
    @("O(log n)", " == ""logarithmic"""),
    @("O(n)", " == ""linear"""),
    @("O(n^2)", " == ""quadratic"""),
    @("O(n^3)", " == ""cubic"""),
    @("O(2^n)", " == ""exponential""")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $prefix = $rows[$i][0]
    $suffix = $rows[$i][1]
    $full = $prefix + $suffix

    $cell = $ws.Cells.Item($r, 1)

    # Base cell font: "Var(--ff-mono)" 13pt, dark grey.
    $cell.Font.Name = "Var(--ff-mono)"
    $cell.Font.Size = 13
    $cell.Font.Color = $darkColor

    $cell.Value = $full

    # Second run (" == ""word""") gets its own run-level formatting.
    $chars = $cell.Characters($prefix.Length + 1, $suffix.Length)
    $chars.Font.Name = "Inherit"
    $chars.Font.Size = 15
    $chars.Font.Color = $darkColor

    $ws.Rows.Item($r).RowHeight = 19
}

# Row 7: hyperlink to the Wikipedia time-complexity article, with the
# actual cell text restored afterwards (the hyperlink's display text
# differs from the literal cell contents in the source workbook).
$ws.Hyperlinks.Add($ws.Cells.Item(7, 1), "http://en.wikipedia.org/wiki/Time_complexity", "Linearithmic.2Fquasilinear_time", $null, "http://en.wikipedia.org/wiki/Time_complexity - Linearithmic.2Fquasilinear_time")
$ws.Cells.Item(7, 1).Value = "O(n log n) == ""linearithmic"""

# Match the saved selection/view state of the sheet.
$ws.Range("D10").Select()
